$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet
$ws.Name = "Weekly_Update"

# 2. Apply number format + centered alignment to columns B:D (entire columns)
$cols = $ws.Range("B:D")
$cols.NumberFormat = "0.00"
$cols.HorizontalAlignment = -4108  # xlCenter

# 3. Add the new "Cash" row beneath the existing data (row 19)
$ws.Range("A19").Value = "Cash"
$ws.Range("B19").Value = "---"
$ws.Range("C19").Value = 1229.532809899639
$ws.Range("D19").Value = 1335.176621862619

$ws.Range("B19:D19").NumberFormat = "0.00"
$ws.Range("B19:D19").HorizontalAlignment = -4108  # xlCenter
